# Apply crypto price/volume updates per Sat Jan 27 20:24:35 UTC 2024 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text figures (e.g. "41.969.89"); force
# text format before assigning so Excel does not reinterpret/round them as numbers.
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextCell "D2" '41.969.89'
Set-TextCell "E2" '  -0.30%  '
Set-TextCell "D3" '2.272.25'
Set-TextCell "E3" '  +0.61%  '
Set-TextCell "E4" '  -0.10%  '
Set-TextCell "D5" '306.10'
Set-TextCell "E5" '  +1.23%  '
Set-TextCell "D6" '93.06'
Set-TextCell "E6" '  -0.11%  '
Set-TextCell "E7" '  -0.45%  '
Set-TextCell "E8" '  -0.14%  '
Set-TextCell "E9" '  +0.86%  '
Set-TextCell "D10" '32.81'
Set-TextCell "E10" '  -0.03%  '
Set-TextCell "D11" '0.0802'
Set-TextCell "E11" '  +0.07%  '
Set-TextCell "E12" '  -2.03%  '
Set-TextCell "D13" '6.69'
Set-TextCell "E13" '  -0.27%  '
Set-TextCell "D14" '2.623.90'
Set-TextCell "E14" '  +0.55%  '
Set-TextCell "E15" '  +1.30%  '
Set-TextCell "D16" '2.272.93'
Set-TextCell "E16" '  -0.22%  '
Set-TextCell "D17" '0.786'
Set-TextCell "E17" '  +3.61%  '
Set-TextCell "D18" '41.883.12'
Set-TextCell "E18" '  -0.25%  '
Set-TextCell "D19" '12.78'
Set-TextCell "E19" '  +4.01%  '
Set-TextCell "D20" '0.0₃0919'
Set-TextCell "E20" '  +1.50%  '
Set-TextCell "E21" '  +0.53%  '
Set-TextCell "D22" '68.24'
Set-TextCell "E22" '  +1.47%  '
Set-TextCell "D23" '244.29'
Set-TextCell "E23" '  +0.91%  '
Set-TextCell "E24" '  +0.10%  '
Set-TextCell "D25" '1.95'
Set-TextCell "E25" '  +0.64%  '
Set-TextCell "E26" '  +0.01%  '
Set-TextCell "E27" '  -0.10%  '
Set-TextCell "D28" '9.71'
Set-TextCell "E28" '  +0.08%  '
Set-TextCell "D29" '2.08'
Set-TextCell "E29" '  -4.43%  '
Set-TextCell "D30" '35.10'
Set-TextCell "E30" '  +2.59%  '
Set-TextCell "D31" '159.14'
Set-TextCell "E31" '  +0.09%  '
Set-TextCell "D32" '5.36'
Set-TextCell "E32" '  +3.82%  '
Set-TextCell "E33" '  -0.06%  '
Set-TextCell "D34" '0.0745'
Set-TextCell "E34" '  -0.17%  '
Set-TextCell "E35" '  -0.43%  '
Set-TextCell "D36" '17.30'
Set-TextCell "E36" '  +3.31%  '
Set-TextCell "E37" '  -1.30%  '
Set-TextCell "E38" '  +0.07%  '
Set-TextCell "D39" '0.117'
Set-TextCell "E39" '  +0.78%  '
Set-TextCell "E40" '  -0.35%  '
Set-TextCell "D41" '3.97'
Set-TextCell "E41" '  +0.60%  '
Set-TextCell "D42" '19.87'
Set-TextCell "E42" '  -0.95%  '
Set-TextCell "D43" '2.015.42'
Set-TextCell "E43" '  -1.78%  '
Set-TextCell "E44" '  +9.31%  '
Set-TextCell "E45" '  +0.92%  '
Set-TextCell "D46" '10.34'
Set-TextCell "E46" '  +2.05%  '
Set-TextCell "D47" '2.93'
Set-TextCell "E47" '  +1.72%  '
Set-TextCell "D48" '53.35'
Set-TextCell "E48" '  +2.60%  '
Set-TextCell "D51" '1.15'
Set-TextCell "E51" '  +0.30%  '

# Rows 49/50: BitcoinSV and Stacks swap rank positions, with refreshed figures
Set-TextCell "B49" 'Stacks'
Set-TextCell "C49" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D49" '1.52'
Set-TextCell "E49" '  -0.86%  '
Set-TextCell "B50" 'BitcoinSV'
Set-TextCell "C50" 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextCell "D50" '72.58'
Set-TextCell "E50" '  +2.99%  '
